$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1)
$pf = $p1.Format

# Add a paragraph border (top/left/bottom/right) that reserves 5pt of space
# on every side (no visible line - matches w:pBdr with only w:space set).
$b = $pf.Borders
$b.DistanceFromTop = 5
$b.DistanceFromBottom = 5
$b.DistanceFromLeft = 5
$b.DistanceFromRight = 5

# Widen the left indent from 120 twips (6pt) to 225 twips (11.25pt)
$pf.LeftIndent = 11.25

# Update the placeholder id text and drop the extra run that only
# contained a trailing space character.
$d.Content.Find.Execute("**ID__AFFARS_5309_topic_7__ID** ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_SUBPART_5309_2__ID**", 2)
